$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A25").Value = "Custom Sphere"
$ws.Range("A25").Font.Bold = $true
$ws.Range("B25").Value = "Just considering F_z"

$ws.Range("A26").Value = 15
$ws.Range("B26").Value = 1.0640000000000001
$ws.Range("C26").Formula = "=((((1.4435*(B26*B26))/((B26*B26)-0.020216))+1)^0.5)/1.328"
$ws.Range("D26").Value = 2
$ws.Range("E26").Value = 0.089647142400000004
$ws.Range("G26").Value = "Final Results for custom Polystyrene Bead in water (Radius 1 micro m)"

$ws.Range("C7").Select()
